$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Theme accent colors: accent1 <-> accent5 swap ---
$wb.Theme.ThemeColorScheme.Colors(5).RGB = 13998939   # accent1 -> 5B9BD5
$wb.Theme.ThemeColorScheme.Colors(9).RGB = 12874308   # accent5 -> 4472C4

# --- Default workbook font: Calibri -> Arial ---
$wb.Styles("Normal").Font.Name = "Arial"

# --- Drop the unused trailing columns (C, D) that held the old headers ---
$ws.Columns("A:D").ClearFormats()
$ws.Columns("C:D").Delete()

# --- New header content replacing the old TerminalId/ErrandTypeId/AssigneeId ---
$ws.Range("A1").Value = "TicketId"
$ws.Range("B1").Value = "Action"

# --- Column widths for the two remaining columns ---
$ws.Columns("A").ColumnWidth = 14.5
$ws.Columns("B").ColumnWidth = 17.5

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Selection moves to B2 ---
$ws.Range("B2").Select() | Out-Null
